$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Helper: write a plain text value into a cell without ever letting Excel
# reinterpret numeric-looking text (e.g. "517.19") as a real number, and
# without leaving any residual number-format/style on the cell. We do this
# by momentarily writing a formula that returns the literal string, then
# converting that formula to its static value via copy / paste-special
# (values only) - exactly as typing the value and using Paste Values would.
function Set-TextValue($ws, $addr, $text) {
    $escaped = $text.Replace('"', '""')
    $ws.Range($addr).Formula = '="' + $escaped + '"'
    $ws.Range($addr).Copy() | Out-Null
    $ws.Range($addr).PasteSpecial(-4163) | Out-Null
}


Set-TextValue $ws 'D2' '60.538.35'
Set-TextValue $ws 'E2' '  +0.18%  '

Set-TextValue $ws 'D3' '2.594.00'
Set-TextValue $ws 'E3' '  +0.20%  '

Set-TextValue $ws 'E4' '  +0.05%  '

Set-TextValue $ws 'D5' '517.19'
Set-TextValue $ws 'E5' '  +1.98%  '

Set-TextValue $ws 'D6' '153.55'
Set-TextValue $ws 'E6' '  -0.14%  '

Set-TextValue $ws 'E7' '  +0.05%  '

Set-TextValue $ws 'E8' '  +3.19%  '

Set-TextValue $ws 'D9' '6.69'
Set-TextValue $ws 'E9' '  +0.68%  '

Set-TextValue $ws 'E10' '  +1.83%  '

Set-TextValue $ws 'E11' '  +0.32%  '

Set-TextValue $ws 'E12' '  +1.39%  '

Set-TextValue $ws 'D13' '3.051.38'
Set-TextValue $ws 'E13' '  +0.34%  '

Set-TextValue $ws 'D14' '60.554.30'
Set-TextValue $ws 'E14' '  +0.31%  '

Set-TextValue $ws 'D15' '21.67'
Set-TextValue $ws 'E15' '  +0.24%  '

Set-TextValue $ws 'E16' '  +0.37%  '

Set-TextValue $ws 'D17' '2.600.94'
Set-TextValue $ws 'E17' '  +0.62%  '

Set-TextValue $ws 'E18' '  -1.42%  '

Set-TextValue $ws 'D19' '350.94'
Set-TextValue $ws 'E19' '  +1.23%  '

Set-TextValue $ws 'D20' '10.56'
Set-TextValue $ws 'E20' '  +1.88%  '

Set-TextValue $ws 'E21' '  +1.50%  '

Set-TextValue $ws 'E22' '  +0.15%  '

Set-TextValue $ws 'D23' '60.94'

Set-TextValue $ws 'D24' '0.427'
Set-TextValue $ws 'E24' '  +1.59%  '

Set-TextValue $ws 'E25' '  +0.14%  '

Set-TextValue $ws 'D26' '2.713.28'
Set-TextValue $ws 'E26' '  +0.44%  '

Set-TextValue $ws 'E27' '  +0.34%  '

Set-TextValue $ws 'D28' '0.0₃0841'
Set-TextValue $ws 'E28' '  -0.47%  '

Set-TextValue $ws 'D29' '7.33'
Set-TextValue $ws 'E29' '  -1.65%  '

Set-TextValue $ws 'E30' '  +0.00%  '

Set-TextValue $ws 'D31' '6.27'
Set-TextValue $ws 'E31' '  +9.07%  '

Set-TextValue $ws 'E32' '  +0.26%  '

Set-TextValue $ws 'E33' '  +2.38%  '

Set-TextValue $ws 'D34' '149.83'
Set-TextValue $ws 'E34' '  -2.89%  '

Set-TextValue $ws 'D35' '4.15'
Set-TextValue $ws 'E35' '  +3.97%  '

Set-TextValue $ws 'D36' '1.19'
Set-TextValue $ws 'E36' '  +0.49%  '

Set-TextValue $ws 'E37' '  +7.95%  '

Set-TextValue $ws 'D38' '1.48'
Set-TextValue $ws 'E38' '  +1.71%  '

Set-TextValue $ws 'E39' '  +0.21%  '

Set-TextValue $ws 'D40' '36.36'
Set-TextValue $ws 'E40' '  +1.45%  '

Set-TextValue $ws 'D41' '0.838'
Set-TextValue $ws 'E41' '  -1.21%  '

Set-TextValue $ws 'D42' '285.99'
Set-TextValue $ws 'E42' '  -3.57%  '

Set-TextValue $ws 'E43' '  +1.58%  '

Set-TextValue $ws 'D44' '0.622'
Set-TextValue $ws 'E44' '  +0.55%  '

Set-TextValue $ws 'E45' '  -0.66%  '

Set-TextValue $ws 'D46' '0.998'
Set-TextValue $ws 'E46' '  +0.10%  '

Set-TextValue $ws 'D47' '19.51'
Set-TextValue $ws 'E47' '  -0.73%  '

Set-TextValue $ws 'B48' 'VeChain'
Set-TextValue $ws 'C48' 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
Set-TextValue $ws 'D48' '0.0236'
Set-TextValue $ws 'E48' '  +1.02%  '

Set-TextValue $ws 'B49' 'RenderToken'
Set-TextValue $ws 'C49' 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
Set-TextValue $ws 'D49' '4.82'
Set-TextValue $ws 'E49' '  -0.98%  '

Set-TextValue $ws 'E50' '  +0.19%  '

Set-TextValue $ws 'D51' '18.96'
Set-TextValue $ws 'E51' '  +7.24%  '

$excel.CutCopyMode = $false
